$d = $word.ActiveDocument

# --- Step 1: P2 -- append a full-width period to the end of the weather text. ---
$d.Content.Find.Execute(
    "多云，今天是六一儿童节，又是开心的一天呢", $true, $false, $false, $false, $false,
    $true, 1, $false, "多云，今天是六一儿童节，又是开心的一天呢。", 2) | Out-Null

# --- Step 2: insert a brand-new paragraph right after P2 for "2022年6月2日星期四". ---
# InsertParagraphAfter() clones the paragraph/run formatting (rFonts hint=eastAsia) from P2,
# which is exactly what the new paragraph needs.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "2022年6月2日星期四"

# --- Step 3: the old "2022n年6月2日星期四" paragraph (now paragraph 4) becomes the
# "rain" paragraph that used to be the last one; the rFonts hint (eastAsia) is unchanged. ---
$d.Content.Find.Execute(
    "2022n年6月2日星期四", $true, $false, $false, $false, $false,
    $true, 1, $false, "中雨，今天是农历五月初四，明天就是端午节了。", 2) | Out-Null

# --- Step 4: insert a new paragraph after that one for "2022年6月3日星期五". ---
# It clones rFonts hint=eastAsia from its predecessor, matching the target.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "2022年6月3日星期五"

# --- Step 5: the original last paragraph (now paragraph 6) keeps its rFonts hint="default"
# pPr untouched, only its text changes to the new Dragon Boat Festival line. Scope the
# search to paragraph 6 specifically -- paragraph 4 now holds identical old text, so a
# document-wide Find would hit the wrong paragraph. ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Find.Execute(
    "中雨，今天是农历五月初四，明天就是端午节了。", $true, $false, $false, $false, $false,
    $true, 1, $false, "中雨，今天是农历五月初五，中国传统端午节。", 2) | Out-Null

# --- Step 6: relocate the _GoBack bookmark from the end of (old) paragraph 4 to the end
# of paragraph 2's text. Adding a zero-width bookmark directly at "end-of-text" triggers a
# runtime quirk, so stash it around a throwaway character and delete the character after,
# which leaves the bookmark collapsed exactly where we need it. ---
$p2 = $d.Paragraphs.Item(2)
$insPos = $p2.Range.End - 1
$placeholder = $d.Range($insPos, $insPos)
$placeholder.InsertAfter("Z")
$bmRange = $d.Range($insPos, $insPos + 1)
$bmRange.Bookmarks.Add("_GoBack")
$delRange = $d.Range($insPos, $insPos + 1)
$delRange.Text = ""

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
